# Updates worksheet cells with refreshed Market Board pricing data
# (currentAveragePrice / LevePrice / LeveProfit columns H-N) for several leves
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 10297.5
$ws.Range("I28").Value = 20415
$ws.Range("K28").Value = 20415
$ws.Range("M28").Value = -19930

$ws.Range("H63").Value = 64998
$ws.Range("J63").Value = 64998
$ws.Range("L63").Value = 64998
$ws.Range("N63").Value = -66246

$ws.Range("H66").Value = 64998
$ws.Range("J66").Value = 64998
$ws.Range("L66").Value = 194994
$ws.Range("N66").Value = -201234

$ws.Range("H88").Value = 4940
$ws.Range("J88").Value = 4940
$ws.Range("L88").Value = 4940
$ws.Range("N88").Value = -5752

$ws.Range("H91").Value = 4940
$ws.Range("J91").Value = 4940
$ws.Range("L91").Value = 4940
$ws.Range("N91").Value = -7748

$ws.Range("H106").Value = 6501263
$ws.Range("I106").Value = 7264470
$ws.Range("K106").Value = 7264470
$ws.Range("M106").Value = -7263839

$ws.Range("H138").Value = 2998.9092
$ws.Range("I138").Value = 992.04346
$ws.Range("J138").Value = 4441.3438
$ws.Range("K138").Value = 2976.13038
$ws.Range("L138").Value = 13324.0314
$ws.Range("M138").Value = 2163.86962
$ws.Range("N138").Value = -23604.0314

$ws.Range("H141").Value = 8949.691999999999
$ws.Range("I141").Value = 10061.833
$ws.Range("J141").Value = 7996.4287
$ws.Range("K141").Value = 30185.499
$ws.Range("L141").Value = 23989.2861
$ws.Range("M141").Value = -25005.499
$ws.Range("N141").Value = -34349.2861

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 4595.143
$ws.Range("I45").Value = 4184
$ws.Range("J45").Value = 5047.4
$ws.Range("K45").Value = 4184
$ws.Range("L45").Value = 5047.4
$ws.Range("M45").Value = -3807
$ws.Range("N45").Value = -5801.4

$ws.Range("H61").Value = 4469.93
$ws.Range("I61").Value = 3844.1936
$ws.Range("K61").Value = 3844.1936
$ws.Range("M61").Value = -3632.1936

$ws.Range("H110").Value = 8023.5
$ws.Range("I110").Value = 9244.5
$ws.Range("J110").Value = 3750
$ws.Range("K110").Value = 9244.5
$ws.Range("L110").Value = 3750
$ws.Range("M110").Value = -7199.5
$ws.Range("N110").Value = -7840

$ws.Range("H132").Value = 3432
$ws.Range("I132").Value = 2045.6923
$ws.Range("K132").Value = 6137.0769
$ws.Range("M132").Value = -3607.0769

$ws.Range("H136").Value = 4469.93
$ws.Range("I136").Value = 3844.1936
$ws.Range("K136").Value = 11532.5808
$ws.Range("M136").Value = -8982.5808

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3065.2144
$ws.Range("J20").Value = 3510.3635
$ws.Range("L20").Value = 3510.3635
$ws.Range("N20").Value = -4004.3635

$ws.Range("H134").Value = 2732.077
$ws.Range("I134").Value = 2287.1428
$ws.Range("K134").Value = 6861.428400000001
$ws.Range("M134").Value = -4326.428400000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2250.75
$ws.Range("I16").Value = 2250.75
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 2250.75
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -1963.75
$ws.Range("N16").Value = $null

$ws.Range("H31").Value = 2769.7778
$ws.Range("I31").Value = 1275.4286
$ws.Range("K31").Value = 1275.4286
$ws.Range("M31").Value = -980.4286

$ws.Range("H34").Value = 2769.7778
$ws.Range("I34").Value = 1275.4286
$ws.Range("K34").Value = 1275.4286
$ws.Range("M34").Value = -1073.4286

$ws.Range("H94").Value = 2485.0908
$ws.Range("J94").Value = 1974
$ws.Range("L94").Value = 1974
$ws.Range("N94").Value = -2876

$ws.Range("H113").Value = 2250.75
$ws.Range("I113").Value = 2250.75
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2250.75
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -80.75
$ws.Range("N113").Value = $null

$ws.Range("H134").Value = 8877.4
$ws.Range("I134").Value = 8142.4287
$ws.Range("K134").Value = 24427.2861
$ws.Range("M134").Value = -21892.2861

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 413.7143
$ws.Range("J11").Value = 266.33334
$ws.Range("L11").Value = 799.0000200000001
$ws.Range("N11").Value = -1079.00002

$ws.Range("H12").Value = 197.35294
$ws.Range("J12").Value = 139.42857
$ws.Range("L12").Value = 418.28571
$ws.Range("N12").Value = -764.28571

$ws.Range("H121").Value = 1668089
$ws.Range("J121").Value = 2001616.8
$ws.Range("L121").Value = 6004850.4
$ws.Range("N121").Value = -6007470.4

$ws.Range("H126").Value = 14789.286

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3178.52
$ws.Range("I132").Value = 3457.762
$ws.Range("J132").Value = 1712.5
$ws.Range("K132").Value = 10373.286
$ws.Range("L132").Value = 5137.5
$ws.Range("M132").Value = -7843.286
$ws.Range("N132").Value = -10197.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 16003
$ws.Range("I3").Value = 19004.5
$ws.Range("J3").Value = 10000
$ws.Range("K3").Value = 19004.5
$ws.Range("L3").Value = 10000
$ws.Range("M3").Value = -18892.5
$ws.Range("N3").Value = -10224

$ws.Range("H15").Value = 16003
$ws.Range("I15").Value = 19004.5
$ws.Range("J15").Value = 10000
$ws.Range("K15").Value = 19004.5
$ws.Range("L15").Value = 10000
$ws.Range("M15").Value = -18834.5
$ws.Range("N15").Value = -10340

$ws.Range("H47").Value = 35000
$ws.Range("I47").Value = 21000
$ws.Range("J47").Value = 49000
$ws.Range("K47").Value = 21000
$ws.Range("L47").Value = 49000
$ws.Range("M47").Value = -20510
$ws.Range("N47").Value = -49980

$ws.Range("H52").Value = 35000
$ws.Range("I52").Value = 21000
$ws.Range("J52").Value = 49000
$ws.Range("K52").Value = 21000
$ws.Range("L52").Value = 49000
$ws.Range("M52").Value = -20767
$ws.Range("N52").Value = -49466

$ws.Range("H122").Value = 3712.9062
$ws.Range("I122").Value = 3275.9546
$ws.Range("J122").Value = 4674.2
$ws.Range("K122").Value = 9827.863799999999
$ws.Range("L122").Value = 14022.6
$ws.Range("M122").Value = -7377.863799999999
$ws.Range("N122").Value = -18922.6

$ws.Range("H132").Value = 406149.8
$ws.Range("I132").Value = 711731.3
$ws.Range("J132").Value = 5074.0625
$ws.Range("K132").Value = 2135193.9
$ws.Range("L132").Value = 15222.1875
$ws.Range("M132").Value = -2132663.9
$ws.Range("N132").Value = -20282.1875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 1630166.6
$ws.Range("J9").Value = 2376250
$ws.Range("L9").Value = 2376250
$ws.Range("N9").Value = -2376530

$ws.Range("H14").Value = 1624.75
$ws.Range("I14").Value = 2750
$ws.Range("J14").Value = 499.5
$ws.Range("K14").Value = 2750
$ws.Range("L14").Value = 499.5
$ws.Range("M14").Value = -2582
$ws.Range("N14").Value = -835.5

$ws.Range("H122").Value = 7372.9287
$ws.Range("I122").Value = 4320.1763
$ws.Range("J122").Value = 12090.818
$ws.Range("K122").Value = 12960.5289
$ws.Range("L122").Value = 36272.454
$ws.Range("M122").Value = -10510.5289
$ws.Range("N122").Value = -41172.454
